# Generate Report for Handoff
#
# The localization status report is regenerated: the file
# "416f48b4-653f-46a1-9368-33089e2614cc.md" moves from
# "Handed back: in sync with en-US" to "Ready for handoff" (a new handoff
# round was started for it), while "a6869899-3547-46eb-ba75-ae8e59a522d9.md"
# keeps its "Handed back: in sync with en-US" status. As a side effect the
# two rows on every sheet swap positions (a6869899 now listed first).
#
# Only the cells whose value actually changes are touched, so untouched
# cells (including blanks) keep their original representation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2 now describes a6869899 (still "Handed back: in sync with en-US")
$ws.Range("A2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.md"
$ws.Range("B2").Value = "e2e\a6869899-3547-46eb-ba75-ae8e59a522d9.md"

# Row 3 now describes 416f48b4, which is now "Ready for handoff"
$ws.Range("A3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.md"
$ws.Range("B3").Value = "e2e\416f48b4-653f-46a1-9368-33089e2614cc.md"
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-09-03 04:53:11"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 now describes a6869899 (still "Handed back: in sync with en-US")
$ws.Range("A2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.md"
$ws.Range("G2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.03e74bf2ceace2fa73293d73198009e554b79884.zh-cn.xlf"
$ws.Range("I2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.md"
$ws.Range("J2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.03e74bf2ceace2fa73293d73198009e554b79884.zh-cn.xlf"

# Row 3 now describes 416f48b4, which is now "Ready for handoff"
$ws.Range("A3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("G3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.624efe40e85c3a0a91593478ec01b79f3ac83a7c.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-03 04:53:05"
$ws.Range("I3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.md"
$ws.Range("J3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.624efe40e85c3a0a91593478ec01b79f3ac83a7c.zh-cn.xlf"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b23aea7a43fc87d83ddeaba54a15383b0763d8e9/e2e/416f48b4-653f-46a1-9368-33089e2614cc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3d0855c0b80110657c3c74a4415c422c90befd/e2e/416f48b4-653f-46a1-9368-33089e2614cc.md."

# Column P widened to fit the new error-detail text
$ws.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 now describes a6869899 (still "Handed back: in sync with en-US")
$ws.Range("A2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.md"
$ws.Range("G2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.03e74bf2ceace2fa73293d73198009e554b79884.de-de.xlf"
$ws.Range("I2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.md"
$ws.Range("J2").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.03e74bf2ceace2fa73293d73198009e554b79884.de-de.xlf"

# Row 3 now describes 416f48b4, which is now "Ready for handoff"
$ws.Range("A3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("G3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.624efe40e85c3a0a91593478ec01b79f3ac83a7c.de-de.xlf"
$ws.Range("H3").Value = "2016-09-03 04:53:11"
$ws.Range("I3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.md"
$ws.Range("J3").Value = "416f48b4-653f-46a1-9368-33089e2614cc.624efe40e85c3a0a91593478ec01b79f3ac83a7c.de-de.xlf"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b23aea7a43fc87d83ddeaba54a15383b0763d8e9/e2e/416f48b4-653f-46a1-9368-33089e2614cc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3d0855c0b80110657c3c74a4415c422c90befd/e2e/416f48b4-653f-46a1-9368-33089e2614cc.md."

# Column P widened to fit the new error-detail text
$ws.Columns.Item(16).ColumnWidth = 39.17
